# Update the "Förändrad" date column (C) from 2023-10-03 (45202) to
# 2023-10-04 (45203) for every data row (rows 2 through 330).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 330
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
